# The accounting sheet "Taul1" lists generated PDF report filenames.
# A new file was produced before the existing "RoboCamp" row (row 20),
# so record it in the previously-empty A1 cell, matching the formatting
# already used for the other filename rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

$ws.Range("A1").Value = "Ohjelmistorobotti.Maaliskuu.pdf"

# Copy the existing filename cell's formatting onto the new cell so it
# matches the style used by A20/A21 (style index 4).
$ws.Range("A20").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
